# Update comparison worksheet with newer biosteam/biorefinery module results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comparison")

# --- Updated raw results (columns C:H) for rows 2-8 ---
# Row 2 (sugarcane1g)
$ws.Range("E2").Value = 1.4245895770530499
$ws.Range("F2").Value = 2.6959830587204801
$ws.Range("G2").Value = 1.44065126995148
$ws.Range("H2").Value = 2.5803373928586102

# Row 3 (oilcane1g)
$ws.Range("C3").Value = 2.0069969392743201
$ws.Range("D3").Value = -1.49834624704022
$ws.Range("E3").Value = 1.99400148127827
$ws.Range("F3").Value = -1.5789338598814799
$ws.Range("G3").Value = 2.02343360509968
$ws.Range("H3").Value = -1.93520940173647

# Row 4 (cornstover)
$ws.Range("C4").Value = 1.62033891808155
$ws.Range("D4").Value = -13.0790883729995
$ws.Range("E4").Value = 1.5984241061734099
$ws.Range("F4").Value = -13.418638499597799
$ws.Range("G4").Value = 1.70276209508825
$ws.Range("H4").Value = -14.065598314481599

# Row 5 (sugarcane2g)
$ws.Range("C5").Value = 2.2043606950162902
$ws.Range("E5").Value = 2.2047743971404299
$ws.Range("F5").Value = 1.42090001234551
$ws.Range("G5").Value = 1.77521719675212
$ws.Range("H5").Value = -2.3757211244144001

# Row 6 (oilcane2g)
$ws.Range("C6").Value = 2.6284154714118699
$ws.Range("D6").Value = 1.82286666134993
$ws.Range("E6").Value = 2.6286053641353702
$ws.Range("F6").Value = 1.90308498431082
$ws.Range("G6").Value = 2.29742812831527
$ws.Range("H6").Value = -0.33666773212129097

# Row 7 (lactic)
$ws.Range("C7").Value = 2.7627905628818201
$ws.Range("D7").Value = 1.5614357145042701
$ws.Range("E7").Value = 2.76279719460446
$ws.Range("F7").Value = 1.6156982112587499
$ws.Range("G7").Value = 2.4922967876054898
$ws.Range("H7").Value = -0.86378104562301405

# Row 8 (average)
$ws.Range("C8").Value = 1.43328539923927
$ws.Range("D8").Value = 4.6067986222162496
$ws.Range("E8").Value = 1.3905358154549301
$ws.Range("F8").Value = 4.6038332682662597
$ws.Range("G8").Value = 1.27394079067161
$ws.Range("H8").Value = 3.3845987314102701

# --- Cursor/selection moved by the author while reviewing the refreshed results ---
$ws.Range("L18").Select()
